$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''42.908.49'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '''2.547.73'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.10%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '''304.01'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.56%  '
$ws.Range("D6").Value = '''98.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.95%  '
$ws.Range("E7").Value = '  +0.75%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").Value = '''0.545'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.47%  '
$ws.Range("D10").Value = '''37.17'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.03%  '
$ws.Range("D11").Value = '''0.0824'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.31%  '
$ws.Range("E12").Value = '  +5.63%  '
$ws.Range("E13").Value = '  -0.99%  '
$ws.Range("D14").Value = '''2.945.42'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.42%  '
$ws.Range("D15").Value = '''2.565.15'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.96%  '
$ws.Range("E16").Value = '  +7.11%  '
$ws.Range("E17").Value = '  +1.04%  '
$ws.Range("D18").Value = '''42.935.93'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("D19").Value = '''13.83'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.71%  '
$ws.Range("D20").Value = '''0.0₃0994'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.47%  '
$ws.Range("D21").Value = '''6.59'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.22%  '
$ws.Range("D22").Value = '''71.94'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.22%  '
$ws.Range("D23").Value = '''254.68'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.20%  '
$ws.Range("E24").Value = '  +1.81%  '
$ws.Range("E25").Value = '  -1.79%  '
$ws.Range("D26").Value = '''28.04'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.43%  '
$ws.Range("D27").Value = '''1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.09%  '
$ws.Range("D28").Value = '''10.25'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.26%  '
$ws.Range("D29").Value = '''37.81'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.65%  '
$ws.Range("D30").Value = '''2.08'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.98%  '
$ws.Range("D31").Value = '''6.18'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.88%  '
$ws.Range("D32").Value = '''158.89'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.43%  '
$ws.Range("D33").Value = '''19.61'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +14.88%  '
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("E35").Value = '  +0.76%  '
$ws.Range("E36").Value = '  -1.95%  '
$ws.Range("E37").Value = '  -4.66%  '
$ws.Range("D38").Value = '''0.117'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.56%  '
$ws.Range("D39").Value = '''25.46'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.59%  '
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("D41").Value = '''2.11'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +32.65%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = '''3.90'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.40%  '
$ws.Range("B43").Value = 'NEARProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D43").Value = '''3.42'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.35%  '
$ws.Range("D44").Value = '''0.0306'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.14%  '
$ws.Range("D45").Value = '''2.091.15'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.65%  '
$ws.Range("D46").Value = '''1.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").Value = '''86.59'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.33%  '
$ws.Range("D48").Value = '''8.98'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.07%  '
$ws.Range("D49").Value = '''2.802.17'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.35%  '
$ws.Range("D50").Value = '''74.91'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.83%  '
$ws.Range("D51").Value = '''103.35'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.09%  '
